$data = @(
  @(16, "CC", "73127237", "ALEX RICARDO VARGAS VILLARREAL", "1904", 22916, 781242),
  @(17, "CC", "73127237", "ALEX RICARDO VARGAS VILLARREAL", "1903", 31249, 781242),
  @(18, "CC", "73127237", "ALEX RICARDO VARGAS VILLARREAL", "1902", 31249, 781242),
  @(19, "CC", "73127237", "ALEX RICARDO VARGAS VILLARREAL", "1901", 31249, 781242),
  @(20, "CC", "73127237", "ALEX RICARDO VARGAS VILLARREAL", "1812", 15625, 781242),
  @(21, "CC", "6617517", "WILLIAM EGUITH MEZA FIGUEROA", "1904", 22916, 689455),
  @(22, "CC", "6617517", "WILLIAM EGUITH MEZA FIGUEROA", "1903", 31249, 689455),
  @(23, "CC", "6617517", "WILLIAM EGUITH MEZA FIGUEROA", "1902", 31249, 689455),
  @(24, "CC", "6617517", "WILLIAM EGUITH MEZA FIGUEROA", "1901", 31249, 689455),
  @(25, "CC", "6617517", "WILLIAM EGUITH MEZA FIGUEROA", "1812", 15625, 689455),
  @(26, "CC", "73184981", "JESUS ENRIQUE ARZUZA BUELVAS", "1904", 22916, 781242),
  @(27, "CC", "73184981", "JESUS ENRIQUE ARZUZA BUELVAS", "1903", 31249, 781242),
  @(28, "CC", "73184981", "JESUS ENRIQUE ARZUZA BUELVAS", "1902", 31249, 781242),
  @(29, "CC", "73184981", "JESUS ENRIQUE ARZUZA BUELVAS", "1901", 31249, 781242),
  @(30, "CC", "73184981", "JESUS ENRIQUE ARZUZA BUELVAS", "1812", 15625, 781242),
  @(31, "CC", "73150538", "WILMER RAFAEL ARZUZA DIAZ", "1904", 22916, 781242),
  @(32, "CC", "73150538", "WILMER RAFAEL ARZUZA DIAZ", "1903", 31249, 781242),
  @(33, "CC", "73150538", "WILMER RAFAEL ARZUZA DIAZ", "1902", 31249, 781242),
  @(34, "CC", "73150538", "WILMER RAFAEL ARZUZA DIAZ", "1901", 31249, 781242),
  @(35, "CC", "73150538", "WILMER RAFAEL ARZUZA DIAZ", "1812", 14583, 781242),
  @(36, "CC", "73095397", "JUAN ARZUZA HERRERA", "1904", 22916, 781242),
  @(37, "CC", "73570974", "RODRIGO ALFONSO PAJOY CASTILLO", "1904", 22916, 781242),
  @(38, "CC", "73570974", "RODRIGO ALFONSO PAJOY CASTILLO", "1903", 31249, 781242),
  @(39, "CC", "73570974", "RODRIGO ALFONSO PAJOY CASTILLO", "1902", 31249, 781242),
  @(40, "CC", "73570974", "RODRIGO ALFONSO PAJOY CASTILLO", "1901", 31249, 781242),
  @(41, "CC", "73570974", "RODRIGO ALFONSO PAJOY CASTILLO", "1812", 15625, 781242),
  @(42, "CC", "8688130", "JAIME MORENO MORENO", "1904", 22916, 781242),
  @(43, "CC", "8688130", "JAIME MORENO MORENO", "1903", 31249, 781242),
  @(44, "CC", "8688130", "JAIME MORENO MORENO", "1902", 31249, 781242),
  @(45, "CC", "8688130", "JAIME MORENO MORENO", "1901", 31249, 781242),
  @(46, "CC", "8688130", "JAIME MORENO MORENO", "1812", 15625, 781242),
  @(47, "CC", "73193047", "LUIS ALBERTO RODRIGUEZ OSPINA", "1904", 22916, 781242),
  @(48, "CC", "73193047", "LUIS ALBERTO RODRIGUEZ OSPINA", "1903", 31249, 781242),
  @(49, "CC", "73193047", "LUIS ALBERTO RODRIGUEZ OSPINA", "1902", 31249, 781242),
  @(50, "CC", "73193047", "LUIS ALBERTO RODRIGUEZ OSPINA", "1901", 31249, 781242),
  @(51, "CC", "73193047", "LUIS ALBERTO RODRIGUEZ OSPINA", "1812", 14583, 781242),
  @(52, "CC", "3800719", "EDWIN POZUELO ARRIETA", "1904", 22916, 781242),
  @(53, "CC", "3800719", "EDWIN POZUELO ARRIETA", "1903", 31249, 781242),
  @(54, "CC", "3800719", "EDWIN POZUELO ARRIETA", "1902", 31249, 781242),
  @(55, "CC", "3800719", "EDWIN POZUELO ARRIETA", "1901", 31249, 781242),
  @(56, "CC", "3800719", "EDWIN POZUELO ARRIETA", "1812", 15625, 781242),
  @(57, "CC", "14208206", "ISIDRO ORTIZ VARON", "1904", 22916, 781242),
  @(58, "CC", "14208206", "ISIDRO ORTIZ VARON", "1903", 31249, 781242),
  @(59, "CC", "14208206", "ISIDRO ORTIZ VARON", "1902", 31249, 781242),
  @(60, "CC", "14208206", "ISIDRO ORTIZ VARON", "1901", 31249, 781242),
  @(61, "CC", "14208206", "ISIDRO ORTIZ VARON", "1812", 14583, 781242)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    $ws.Cells.Item($r, 4).Value = $item[3]
    $ws.Cells.Item($r, 5).Value = $item[4]
    $ws.Cells.Item($r, 6).Value = $item[5]
    $ws.Cells.Item($r, 7).Value = $item[6]
}
